# Applies the "User Stories" content edit described by the commit diff:
#  - Rephrases several German user-story texts (shared strings) into
#    consistent "Ich-Form" (first-person) wording.
#  - Renames the stakeholder role "Veranstalter" -> "Anbieter" everywhere
#    it appears (both as its own cell and inline inside a longer text).
#  - Updates the StoryPoint / Estimate (hours) numbers in columns K and L
#    for rows 4-16, and clears L15 (now blank).
#  - Renames the built-in cell style from "Standard" to "Normal".
#  - Updates the active selection on the sheet to K15.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Text / wording edits (these live in the shared-string table; setting
#    .Value on the cell that uses each string updates the visible text).
# ---------------------------------------------------------------------
$ws.Range("J4").Value  = "ich sofort Aktivitäten in meiner Nähe sehe"
$ws.Range("M4").Value  = "- App fragt nach Standortberechtigung`n- Standort wird korrekt erkannt`n- Bei Ablehnung kann Standort manuell eingegeben werden"

$ws.Range("J5").Value  = "ich passende Vorschläge in meiner Umgebung finde"
$ws.Range("M5").Value  = "- Radius auswählbar (z.B. 1km, 5km, 10km)`n- Liste aktualisiert sich automatisch`n- Entfernung wird pro Aktivität angezeigt"

$ws.Range("I6").Value  = "nach Kategorien (z.B. Essen, Sport, Kultur, Outdoor) filtern"
$ws.Range("J6").Value  = "ich schneller passende Aktivitäten finde"

$ws.Range("J7").Value  = "ich entscheiden kann, ob sie für mich geeignet ist"

# "Veranstalter" -> "Anbieter" (shared by H8, H9, H10, H11, H14)
$ws.Range("H8").Value  = "Anbieter"
$ws.Range("H9").Value  = "Anbieter"
$ws.Range("H10").Value = "Anbieter"
$ws.Range("H11").Value = "Anbieter"
$ws.Range("H14").Value = "Anbieter"

$ws.Range("J8").Value  = "ich meine Aktivitäten eintragen und verwalten kann"

$ws.Range("I9").Value  = "eine neue Aktivität oder mein Lokal eintragen"
$ws.Range("J9").Value  = "Nutzer mein Angebot in der App finden können"

$ws.Range("I10").Value = "meine Aktivität oder die Informationen zu meinem Geschäft bearbeiten"
$ws.Range("J10").Value = "ich Änderungen aktualisieren kann"

$ws.Range("J11").Value = "sie nicht mehr öffentlich sichtbar ist"

$ws.Range("J12").Value = "sodass ich Aktivitäten im Vorraus planen kann Bewertungen abgeben und Favoriten markieren kann"

$ws.Range("M14").Value = "- Aktivität ist bis zur Endzeit im Status „Aktiv“ und wird Nutzern normal angezeigt`n- Nach der Endzeit wechselt der Status automatisch auf „Abgelaufen/Deaktiviert“                                                                                              -Anbieter sieht abgelaufene Aktivitäten weiterhin unter „Meine Aktivitäten“ (mit Status)"

$ws.Range("J16").Value = "ich nichts vergesse"

# ---------------------------------------------------------------------
# 2. StoryPoint (K) / Estimates-hours (L) number updates, rows 4-16.
# ---------------------------------------------------------------------
$ws.Range("K4").Value  = 8
$ws.Range("L4").Value  = 10

$ws.Range("K5").Value  = 8
$ws.Range("L5").Value  = 10

$ws.Range("K6").Value  = 5
$ws.Range("L6").Value  = 4

$ws.Range("K7").Value  = 5
$ws.Range("L7").Value  = 4

$ws.Range("K8").Value  = 8
$ws.Range("L8").Value  = 10

$ws.Range("K9").Value  = 8
$ws.Range("L9").Value  = 10

$ws.Range("K10").Value = 5
$ws.Range("L10").Value = 4

$ws.Range("K11").Value = 5
$ws.Range("L11").Value = 4

$ws.Range("K12").Value = 2
$ws.Range("L12").Value = 2

$ws.Range("K13").Value = 2
$ws.Range("L13").Value = 2

$ws.Range("K14").Value = 3
$ws.Range("L14").Value = 3

$ws.Range("K15").Value = 2
$ws.Range("L15").ClearContents()

$ws.Range("K16").Value = 1
$ws.Range("L16").Value = 5

# ---------------------------------------------------------------------
# 3. Row 17 height (cosmetic re-wrap after the edits above).
# ---------------------------------------------------------------------
$ws.Rows.Item(17).RowHeight = 106.5

# ---------------------------------------------------------------------
# 4. Built-in cell style rename: "Standard" -> "Normal".
# ---------------------------------------------------------------------
$wb.Styles.Item(1).Name = "Normal"

# ---------------------------------------------------------------------
# 5. Update the active selection to K15 (matches the saved view state).
# ---------------------------------------------------------------------
$ws.Activate()
$ws.Range("K15").Select()
